# Add translation settings (display.locale.text / pt / sw) to the
# "settings" sheet of the ODK-X form, and make "settings" the active
# sheet (it was "model" before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")
$wsm = $wb.Worksheets.Item("model")

# --- Header row: new translation-related setting columns (D1:H1) ---
$ws.Range("D1").Value = "display.title.text.pt"
$ws.Range("E1").Value = "display.title.text.sw"
$ws.Range("F1").Value = "display.locale.text"
$ws.Range("G1").Value = "display.locale.text.pt"
$ws.Range("H1").Value = "display.locale.text.sw"

# --- Row 5 (survey / "Household Members New" title): repeat the same
# display title text for the Portuguese and Swahili columns too ---
$ws.Range("D5").Value = "Household Members New"
$ws.Range("E5").Value = "Household Members New"

# --- New rows 7-9: locale rows (default / pt / sw) ---
$ws.Range("A7").Value = "default"
$ws.Range("F7").Value = "English"
$ws.Range("G7").Value = "English"
$ws.Range("H7").Value = "English"

$ws.Range("A8").Value = "pt"
$ws.Range("F8").Value = "Português"
$ws.Range("G8").Value = "Português"
$ws.Range("H8").Value = "Português"

$ws.Range("A9").Value = "sw"
$ws.Range("F9").Value = "Kiswahili"
$ws.Range("G9").Value = "Kiswahili"
$ws.Range("H9").Value = "Kiswahili"

# --- Formatting ---
# C1 (already present) and the brand-new header/locale cells should use
# the sheet's plain default style (no custom font), same as e.g.
# model!A1. Use a formats-only paste so we reuse the existing "Normal"
# style instead of registering a new one.
$wsm.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("D1:H1").PasteSpecial(-4122)
$ws.Range("A7:A9").PasteSpecial(-4122)
$ws.Range("F7:H9").PasteSpecial(-4122)

# D5/E5 should keep the same font style as C5 (the "Household Members
# New" title cell they duplicate).
$ws.Range("C5").Copy()
$ws.Range("D5:E5").PasteSpecial(-4122)

# --- Column widths for the new columns ---
$ws.Columns.Item(4).ColumnWidth = 20.94
$ws.Columns.Item(5).ColumnWidth = 20.94
$ws.Columns.Item(6).ColumnWidth = 15.5
$ws.Columns.Item(7).ColumnWidth = 17.65
$ws.Columns.Item(8).ColumnWidth = 18.4

# --- Make "settings" the active/selected sheet (was "model") ---
$ws.Activate()
